# Generate Report for Archive
#
# 1) Update the localization status text "Ready for handoff" -> "In Translation"
#    everywhere it is used: Overview!E2, Overview!F2, zh-cn!C2, de-de!C2
# 2) Narrow the "zh-cn"/"de-de" status columns (Overview E:F, zh-cn C, de-de C)
#    from their current width down to the new, narrower width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1) Status text update -------------------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2) Column width update --------------------------------------------------
# Target stored width ~= 13.4101848602295 characters (down from ~17.2159881591797)
$newColumnWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth  # column E
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth  # column F
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth      # column C
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth      # column C
